# Add a new "UID" column (K) with sequential numeric IDs starting at 4001.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1 mirrors the style used by J1 ("Team") -> style index 2.
$ws.Range("K1").Value = "UID"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill K2:K89 with sequential UID values 4001..4088.
$startUid = 4001
for ($row = 2; $row -le 89; $row++) {
    $ws.Cells.Item($row, 11).Value = $startUid + ($row - 2)
}

# Update the view: scroll so column B is the left-most visible column,
# and select the new UID column's data range.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("K2:K89").Select()

# Page setup: A4, portrait.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
